# Natmi following Dr Hou advice:
# Re-run of the LR-pair analysis (Hsp90aa1-Fgfr3) adds "FAPs" as a third
# target cluster (alongside ECs and sCs) and recomputes every expression /
# specificity statistic for all 5x3 sending/target cluster combinations,
# growing the sheet from 10 data rows (A2:T11) to 15 data rows (A2:T16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 15,20
$data[0,0] = "ECs"
$data[0,1] = "Hsp90aa1"
$data[0,2] = "Fgfr3"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 83.148033
$data[0,7] = 249.444099
$data[0,8] = 0.1291567785869119
$data[0,9] = 0.1291567785869119
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 5.243417666666667
$data[0,13] = 15.730253
$data[0,14] = 0.8253998362974575
$data[0,15] = 0.8253998362974574
$data[0,16] = 435.979865180783
$data[0,17] = 3923.818786627047
$data[0,18] = 0.106605983902344
$data[0,19] = 0.106605983902344
$data[1,0] = "ECs"
$data[1,1] = "Hsp90aa1"
$data[1,2] = "Fgfr3"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 83.148033
$data[1,7] = 249.444099
$data[1,8] = 0.1291567785869119
$data[1,9] = 0.1291567785869119
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.6792986666666666
$data[1,13] = 2.037896
$data[1,14] = 0.1069327381315001
$data[1,15] = 0.1069327381315001
$data[1,16] = 56.482347952856
$data[1,17] = 508.341131575704
$data[1,18] = 0.01381108798254238
$data[1,19] = 0.01381108798254239
$data[2,0] = "ECs"
$data[2,1] = "Hsp90aa1"
$data[2,2] = "Fgfr3"
$data[2,3] = "sCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 83.148033
$data[2,7] = 249.444099
$data[2,8] = 0.1291567785869119
$data[2,9] = 0.1291567785869119
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.4298626666666667
$data[2,13] = 1.289588
$data[2,14] = 0.06766742557104236
$data[2,15] = 0.06766742557104236
$data[2,16] = 35.74223519346801
$data[2,17] = 321.680116741212
$data[2,18] = 0.008739706702025457
$data[2,19] = 0.008739706702025458
$data[3,0] = "FAPs"
$data[3,1] = "Hsp90aa1"
$data[3,2] = "Fgfr3"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 154.0559436666667
$data[3,7] = 462.167831
$data[3,8] = 0.2393005425173851
$data[3,9] = 0.2393005425173851
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 5.243417666666667
$data[3,13] = 15.730253
$data[3,14] = 0.8253998362974575
$data[3,15] = 0.8253998362974574
$data[3,16] = 807.7796566768048
$data[3,17] = 7270.016910091244
$data[3,18] = 0.1975186286197424
$data[3,19] = 0.1975186286197424
$data[4,0] = "FAPs"
$data[4,1] = "Hsp90aa1"
$data[4,2] = "Fgfr3"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 154.0559436666667
$data[4,7] = 462.167831
$data[4,8] = 0.2393005425173851
$data[4,9] = 0.2393005425173851
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.6792986666666666
$data[4,13] = 2.037896
$data[4,14] = 0.1069327381315001
$data[4,15] = 0.1069327381315001
$data[4,16] = 104.6499971248418
$data[4,17] = 941.8499741235761
$data[4,18] = 0.02558906224773743
$data[4,19] = 0.02558906224773744
$data[5,0] = "FAPs"
$data[5,1] = "Hsp90aa1"
$data[5,2] = "Fgfr3"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 154.0559436666667
$data[5,7] = 462.167831
$data[5,8] = 0.2393005425173851
$data[5,9] = 0.2393005425173851
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.4298626666666667
$data[5,13] = 1.289588
$data[5,14] = 0.06766742557104236
$data[5,15] = 0.06766742557104236
$data[5,16] = 66.22289876040313
$data[5,17] = 596.0060888436282
$data[5,18] = 0.01619285164990521
$data[5,19] = 0.01619285164990521
$data[6,0] = "M1"
$data[6,1] = "Hsp90aa1"
$data[6,2] = "Fgfr3"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 154.1469826666666
$data[6,7] = 462.4409479999999
$data[6,8] = 0.2394419566139251
$data[6,9] = 0.2394419566139251
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 5.243417666666667
$data[6,13] = 15.730253
$data[6,14] = 0.8253998362974575
$data[6,15] = 0.8253998362974574
$data[6,16] = 808.2570121777603
$data[6,17] = 7274.313109599842
$data[6,18] = 0.1976353517918767
$data[6,19] = 0.1976353517918766
$data[7,0] = "M1"
$data[7,1] = "Hsp90aa1"
$data[7,2] = "Fgfr3"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 154.1469826666666
$data[7,7] = 462.4409479999999
$data[7,8] = 0.2394419566139251
$data[7,9] = 0.2394419566139251
$data[7,10] = 2
$data[7,11] = 0.6666666666666666
$data[7,12] = 0.6792986666666666
$data[7,13] = 2.037896
$data[7,14] = 0.1069327381315001
$data[7,15] = 0.1069327381315001
$data[7,16] = 104.7118397961564
$data[7,17] = 942.4065581654079
$data[7,18] = 0.02560418404429085
$data[7,19] = 0.02560418404429085
$data[8,0] = "M1"
$data[8,1] = "Hsp90aa1"
$data[8,2] = "Fgfr3"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 154.1469826666666
$data[8,7] = 462.4409479999999
$data[8,8] = 0.2394419566139251
$data[8,9] = 0.2394419566139251
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.4298626666666667
$data[8,13] = 1.289588
$data[8,14] = 0.06766742557104236
$data[8,15] = 0.06766742557104236
$data[8,16] = 66.26203302771378
$data[8,17] = 596.358297249424
$data[8,18] = 0.01620242077775753
$data[8,19] = 0.01620242077775753
$data[9,0] = "M2"
$data[9,1] = "Hsp90aa1"
$data[9,2] = "Fgfr3"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 132.8624013333333
$data[9,7] = 398.587204
$data[9,8] = 0.2063798641097711
$data[9,9] = 0.2063798641097711
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 5.243417666666667
$data[9,13] = 15.730253
$data[9,14] = 0.8253998362974575
$data[9,15] = 0.8253998362974574
$data[9,16] = 696.6530623869569
$data[9,17] = 6269.877561482613
$data[9,18] = 0.1703459060512966
$data[9,19] = 0.1703459060512966
$data[10,0] = "M2"
$data[10,1] = "Hsp90aa1"
$data[10,2] = "Fgfr3"
$data[10,3] = "FAPs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 132.8624013333333
$data[10,7] = 398.587204
$data[10,8] = 0.2063798641097711
$data[10,9] = 0.2063798641097711
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.6792986666666666
$data[10,13] = 2.037896
$data[10,14] = 0.1069327381315001
$data[10,15] = 0.1069327381315001
$data[10,16] = 90.25325207586488
$data[10,17] = 812.279268682784
$data[10,18] = 0.02206876396446472
$data[10,19] = 0.02206876396446472
$data[11,0] = "M2"
$data[11,1] = "Hsp90aa1"
$data[11,2] = "Fgfr3"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 132.8624013333333
$data[11,7] = 398.587204
$data[11,8] = 0.2063798641097711
$data[11,9] = 0.2063798641097711
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.4298626666666667
$data[11,13] = 1.289588
$data[11,14] = 0.06766742557104236
$data[11,15] = 0.06766742557104236
$data[11,16] = 57.11258613688356
$data[11,17] = 514.0132752319521
$data[11,18] = 0.01396519409400977
$data[11,19] = 0.01396519409400977
$data[12,0] = "sCs"
$data[12,1] = "Hsp90aa1"
$data[12,2] = "Fgfr3"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 119.5626293333333
$data[12,7] = 358.687888
$data[12,8] = 0.1857208581720069
$data[12,9] = 0.185720858172007
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 5.243417666666667
$data[12,13] = 15.730253
$data[12,14] = 0.8253998362974575
$data[12,15] = 0.8253998362974574
$data[12,16] = 626.9168029195183
$data[12,17] = 5642.251226275665
$data[12,18] = 0.1532939659321979
$data[12,19] = 0.1532939659321979
$data[13,0] = "sCs"
$data[13,1] = "Hsp90aa1"
$data[13,2] = "Fgfr3"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 119.5626293333333
$data[13,7] = 358.687888
$data[13,8] = 0.1857208581720069
$data[13,9] = 0.185720858172007
$data[13,10] = 2
$data[13,11] = 0.6666666666666666
$data[13,12] = 0.6792986666666666
$data[13,13] = 2.037896
$data[13,14] = 0.1069327381315001
$data[13,15] = 0.1069327381315001
$data[13,16] = 81.21873468929422
$data[13,17] = 730.968612203648
$data[13,18] = 0.01985963989246468
$data[13,19] = 0.01985963989246468
$data[14,0] = "sCs"
$data[14,1] = "Hsp90aa1"
$data[14,2] = "Fgfr3"
$data[14,3] = "sCs"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 119.5626293333333
$data[14,7] = 358.687888
$data[14,8] = 0.1857208581720069
$data[14,9] = 0.185720858172007
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.4298626666666667
$data[14,13] = 1.289588
$data[14,14] = 0.06766742557104236
$data[14,15] = 0.06766742557104236
$data[14,16] = 51.39551067890491
$data[14,17] = 462.5595961101441
$data[14,18] = 0.01256725234734439
$data[14,19] = 0.01256725234734439
$ws.Range("A2:T16").Value = $data
